$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "44.519.90"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  +2.57%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.367.33"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.674"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +3.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "238.89"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +2.73%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.25"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +6.42%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +19.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.102"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +7.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "29.35"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +10.76%  "
$ws.Range("E12").Value = "  +1.96%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.724.82"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.85"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +7.94%  "
$ws.Range("E15").Value = "  +7.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.899"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +6.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.373.04"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "44.424.59"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +2.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000103"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +4.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "77.66"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +5.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.46"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +3.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "255.43"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +2.45%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.76"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -3.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.52"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +2.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.46"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +4.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.29"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +0.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.51"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  +0.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "173.87"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -0.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.59"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +5.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.132"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +2.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.132"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +4.81%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0739"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +6.46%  "
$ws.Range("E34").Value = "  +4.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.20"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +3.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.93"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +7.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.43"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -3.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.52"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +0.07%  "
$ws.Range("E39").Value = "  +7.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.82"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +10.31%  "
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.87"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -1.49%  "
$ws.Range("E43").Value = "  +3.85%  "
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0981"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +2.81%  "
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.17"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +0.32%  "
$ws.Range("E46").Value = "  +12.18%  "
$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.47"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +2.00%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "98.59"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -0.80%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.35"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +3.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.442.64"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.593.59"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -0.17%  "
